$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the run-date serial number in M1 (45959 -> 45960)
$ws.Range("M1").Value2 = 45960

# Swap columns B, D, E, F, G between paired rows (these pairs represent
# the same stock item listed twice; the two records were transposed).
$pairs = @(
    @(149,150),
    @(279,280),
    @(313,314),
    @(346,347),
    @(350,352),
    @(372,373),
    @(379,380),
    @(382,383),
    @(419,420),
    @(421,422),
    @(457,458),
    @(536,537),
    @(579,580),
    @(586,587),
    @(593,594),
    @(601,602),
    @(709,710),
    @(720,721),
    @(859,860),
    @(889,890)
)

$cols = @("B","D","E","F","G")

foreach ($pair in $pairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]
    foreach ($col in $cols) {
        $cellA = $ws.Range($col + $rowA)
        $cellB = $ws.Range($col + $rowB)
        $valA = $cellA.Value2
        $valB = $cellB.Value2
        $cellA.Value2 = $valB
        $cellB.Value2 = $valA
    }
}
